$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for the refreshed crypto data.
# NumberFormat '@' + Style reset keeps numeric-looking price strings
# (e.g. '227.57') stored as text, matching the source feed's string cells.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "32.897.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +9.86%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.760.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +6.03%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.69%  "

$ws.Range("E6").Value = "  +4.41%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.66"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +9.23%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.86"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.32%  "

$ws.Range("E10").Value = "  +5.74%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0668"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.50%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0920"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.92%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.010.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.85%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.769.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.43%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.633"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.95%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "10.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.97%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +9.34%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "32.819.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +9.47%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "68.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "259.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.52%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0746"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.89%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.78%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.75%  "

$ws.Range("E25").Value = "  -0.91%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.95%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.66%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.115"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.42%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.64%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.00%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +12.82%  "

$ws.Range("E32").Value = "  +3.43%  "

$ws.Range("E33").Value = "  +5.62%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.57%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.560.20"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.83%  "

$ws.Range("E36").Value = "  +4.50%  "

$ws.Range("E37").Value = "  +2.18%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.633"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.68%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "84.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.36%  "

$ws.Range("E40").Value = "  +5.74%  "

$ws.Range("E41").Value = "  +3.69%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.879"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.54%  "

$ws.Range("E44").Value = "  +6.49%  "

$ws.Range("E45").Value = "  +2.76%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.51%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.87%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.909.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.83%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.10%  "

$ws.Range("E50").Value = "  +5.90%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "95.97"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.66%  "
